# Applies calibration updates to the Honduras model input variables workbook.
# - Updates sampling ranges (max_35 / min_35 columns H / I) for many variable rows
#   across all four strategy sheets.
# - Updates the per-period uncertainty factors (columns J:AS) for the hydropower
#   CO2-capture rows (122-126) on sheet "strategy_id-0" to reflect the new
#   climate-change factor (0.9) that was added.

$wb = $excel.ActiveWorkbook

# Column indices: H=8, I=9, J=10, AS=45
$colH = 8
$colI = 9
$colJ = 10
$colAS = 45

function Set-HI {
    param(
        $ws,
        [int[]]$Rows,
        [double]$HVal,
        [double]$IVal
    )
    foreach ($r in $Rows) {
        $ws.Cells.Item($r, $colH).Value = $HVal
        $ws.Cells.Item($r, $colI).Value = $IVal
    }
}

function Set-JtoAS {
    param(
        $ws,
        [int[]]$Rows,
        [double]$Val
    )
    foreach ($r in $Rows) {
        for ($c = $colJ; $c -le $colAS; $c++) {
            $ws.Cells.Item($r, $c).Value = $Val
        }
    }
}

# ---------------------------------------------------------------------------
# Sheet "strategy_id-0"
# ---------------------------------------------------------------------------
$ws0 = $wb.Worksheets.Item("strategy_id-0")

# Rows 4-23: max_35 1 -> 1.5, min_35 1 -> 0.5
$rows0_15 = 4..23
Set-HI -ws $ws0 -Rows $rows0_15 -HVal 1.5 -IVal 0.5

# Rows where max_35/min_35 converge to 1 / 1
$rows0_11 = @(
    24, 26, 27, 28, 30, 31, 34, 35, 36, 37, 38, 39, 40, 42, 44, 45, 46, 49, 50, 51,
    54, 55, 57, 59, 61, 65, 66, 69, 70, 71, 72, 73, 76, 77, 78, 79, 80, 81, 82, 83,
    84, 85, 86, 87, 88, 89, 90, 92, 93, 94, 95, 96, 97, 98, 99, 100, 101, 102, 104,
    105, 106, 107, 108, 109, 110, 111, 113, 114, 127, 128, 129, 130, 131,
    150, 151, 152, 153, 154, 155, 156, 157, 158, 159, 160, 161, 162, 163, 164,
    169, 170, 171, 172, 173, 174, 175, 176, 177, 178, 179, 180, 181, 182, 183, 184,
    185, 186
)
Set-HI -ws $ws0 -Rows $rows0_11 -HVal 1 -IVal 1

# Rows 122-126 (hydropower CO2-capture): climate change factor added to J:AS -> 0.9
$rows0_jas = 122..126
Set-JtoAS -ws $ws0 -Rows $rows0_jas -Val 0.9

# ---------------------------------------------------------------------------
# Sheet "strategy_id-5001"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("strategy_id-5001")
Set-HI -ws $ws1 -Rows @(2) -HVal 1 -IVal 1

# ---------------------------------------------------------------------------
# Sheet "strategy_id-5006"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("strategy_id-5006")
Set-HI -ws $ws2 -Rows (2..21) -HVal 1.5 -IVal 0.5
Set-HI -ws $ws2 -Rows (22..37) -HVal 1 -IVal 1

# ---------------------------------------------------------------------------
# Sheet "strategy_id-5008"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("strategy_id-5008")
Set-HI -ws $ws3 -Rows (2..21) -HVal 1.5 -IVal 0.5
Set-HI -ws $ws3 -Rows (22..38) -HVal 1 -IVal 1

Write-Host "Edit complete."
